$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 6
$ws.Cells.Item($row, 1).Value = 42608.89334490741
$ws.Cells.Item($row, 2).Value = 18
$ws.Cells.Item($row, 3).Value = 57
$ws.Cells.Item($row, 4).Value = 40
$ws.Cells.Item($row, 5).Value = 73
$ws.Cells.Item($row, 6).Value = 26
$ws.Cells.Item($row, 7).Value = 14394
$ws.Cells.Item($row, 8).Value = 26692
$ws.Cells.Item($row, 9).Value = 3189
$ws.Cells.Item($row, 10).Value = 358
$ws.Cells.Item($row, 11).Value = 255
$ws.Cells.Item($row, 12).Value = 17
$ws.Cells.Item($row, 13).Value = 6
$ws.Cells.Item($row, 14).Value = "Bag"
